$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CleanChassis")

$ws.Range("A2").Value2 = "10123916"
$ws.Range("B2").Value2 = "19XFC2670GE451996"
$ws.Range("A3").Value2 = "10175969"
$ws.Range("B3").Value2 = "1B3CB2HA4AD654852"
$ws.Range("A4").Value2 = "10343518"
$ws.Range("B4").Value2 = "1C4BJWAG3GL281207"
$ws.Range("A5").Value2 = "10078530"
$ws.Range("B5").Value2 = "1C4BJWAG6GL275529"
$ws.Range("A6").Value2 = "10154423"
$ws.Range("B6").Value2 = "1C4BJWAG6HL637523"
$ws.Range("A7").Value2 = "10107943"
$ws.Range("B7").Value2 = "1C4BJWAG7EL135955"
$ws.Range("A8").Value2 = "10205926"
$ws.Range("B8").Value2 = "1C4BJWAG7GL312099"
$ws.Range("A9").Value2 = "10455470"
$ws.Range("B9").Value2 = "1C4BJWBG5CL140020"
$ws.Range("A10").Value2 = "10343825"
$ws.Range("B10").Value2 = "1C4BJWAG7HL659675"
$ws.Range("A11").Value2 = "10238585"
$ws.Range("B11").Value2 = "1C4BJWAG9DL541721"
$ws.Range("A12").Value2 = "10228361"
$ws.Range("B12").Value2 = "1C4BJWAG9DL605076"
$ws.Range("A13").Value2 = "10144267"
$ws.Range("B13").Value2 = "1C4BJWBG1GL123978"
$ws.Range("A14").Value2 = "10527296"
$ws.Range("B14").Value2 = "1C4BJWBGXEL230315"
$ws.Range("A15").Value2 = "10268377"
$ws.Range("B15").Value2 = "1C4BJWBG8FL540982"
$ws.Range("A16").Value2 = "10152681"
$ws.Range("B16").Value2 = "1C4BJWDG6GL149036"
$ws.Range("A17").Value2 = "10254224"
$ws.Range("B17").Value2 = "1C4BJWDG5EL153155"
$ws.Range("A18").Value2 = "10242915"
$ws.Range("B18").Value2 = "1C4BJWDG9JL801655"
$ws.Range("A19").Value2 = "10249326"
$ws.Range("B19").Value2 = "1C4BJWDGXFL734555"
$ws.Range("A20").Value2 = "10043482"
$ws.Range("B20").Value2 = "1C4BJWEG2FL582057"

# Remove the 10 trailing blank rows (509-518) that were part of the
# pre-formatted CleanChassis table range, shrinking it from A1:C518 to A1:C508.
$ws.Range("A509:A518").EntireRow.Delete()
